# Generate Report for Handoff
# Appends a new "ready for handoff" row (row 3) to the Overview, zh-cn and
# de-de sheets/tables, mirroring the existing row 2 for a new source file
# (28ba2baa-...) that is now ready for handoff (was previously "In
# Translation" for the 7635fea2-... file).

$wb = $excel.ActiveWorkbook

$fileNameMd    = "28ba2baa-4d50-416d-8955-9a497193edafooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$pathMd        = "e2e\28ba2baa-4d50-416d-8955-9a497193edafooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$statusText    = "Ready for handoff"
$hoDate        = "2016-08-18 00:26:37"
$zhXlf         = "28ba2baa-4d50-416d-8955-9a497193edafoooooooooooooooooooooooooooooooooooooooo.5c91fa9db51fed0f310ed226352d403ee72ba185.zh-cn.xlf"
$zhXlfDate     = "2016-08-18 00:26:32"
$deXlf         = "28ba2baa-4d50-416d-8955-9a497193edafoooooooooooooooooooooooooooooooooooooooo.5c91fa9db51fed0f310ed226352d403ee72ba185.de-de.xlf"
$hyperlinkUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d7e2e7687a729aff9f6f8af09fa3aaa563dcd498/e2e/28ba2baa-4d50-416d-8955-9a497193edafooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = $fileNameMd
$wsOverview.Range("B3").Value = $pathMd
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText
$wsOverview.Range("G3").Value = $hoDate

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $hyperlinkUrl, "", "", $pathMd) | Out-Null
$wsOverview.Range("B3").Font.Underline = $true
$wsOverview.Range("B3").Font.Color = 13595136

$wsOverview.Columns.Item(5).ColumnWidth = 17
$wsOverview.Columns.Item(6).ColumnWidth = 17

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Range("A3").Value = $fileNameMd
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = $statusText
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "False"
$wsZh.Range("G3").Value = $zhXlf
$wsZh.Range("H3").Value = $zhXlfDate
$wsZh.Range("I3").Value = ""
$wsZh.Range("J3").Value = ""
$wsZh.Range("K3").Value = "0001-01-01 00:00:00"
$wsZh.Range("L3").Value = ""
$wsZh.Range("M3").Value = "True"
$wsZh.Range("N3").Value = ""
$wsZh.Range("O3").Value = "False"
$wsZh.Range("P3").Value = ""

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $hyperlinkUrl, "", "", $fileNameMd) | Out-Null
$wsZh.Range("A3").Font.Underline = $true
$wsZh.Range("A3").Font.Color = 13595136

$wsZh.Columns.Item(3).ColumnWidth = 17

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Range("A3").Value = $fileNameMd
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = $statusText
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "False"
$wsDe.Range("G3").Value = $deXlf
$wsDe.Range("H3").Value = $hoDate
$wsDe.Range("I3").Value = ""
$wsDe.Range("J3").Value = ""
$wsDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDe.Range("L3").Value = ""
$wsDe.Range("M3").Value = "True"
$wsDe.Range("N3").Value = ""
$wsDe.Range("O3").Value = "False"
$wsDe.Range("P3").Value = ""

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $hyperlinkUrl, "", "", $fileNameMd) | Out-Null
$wsDe.Range("A3").Font.Underline = $true
$wsDe.Range("A3").Font.Color = 13595136

$wsDe.Columns.Item(3).ColumnWidth = 17
